# Update countries & provincias Spain
# - Update "last updated" timestamp on the summary cell (A1)
# - Update Suiza (row 18) figures
# - Update Polonia (row 34) figures
# - Swap El Salvador / Jamaica rows (125/126) and refresh their figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 09:22"

# Suiza (row 18)
$ws.Range("E18").Value = 6345
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 1551

# Polonia (row 34)
$ws.Range("D34").Value = 1944
$ws.Range("E34").Value = 8113

# El Salvador now sorts above Jamaica (row 125 <-> row 126), each with
# refreshed case counts.
$ws.Range("A125").Value = "El Salvador"
$ws.Range("B125").Value = 261
$ws.Range("C125").Value = 11
$ws.Range("D125").Value = 67
$ws.Range("E125").Value = 186
$ws.Range("F125").Value = 2
$ws.Range("H125").Value = 8

$ws.Range("A126").Value = "Jamaica"
$ws.Range("B126").Value = 257
$ws.Range("C126").Value = 5
$ws.Range("D126").Value = 28
$ws.Range("E126").Value = 223
$ws.Range("F126").Value = 0
$ws.Range("H126").Value = 6
